# REFAC: merubah format export excel
# Insert a new "No Scan" column before "Tgl Checkin", rename "Total Packing"
# to "Quantity", and update the row data accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "Kode Rak"
$ws.Range("B1").Value = "Tipe Rak"
$ws.Range("C1").Value = "Status Rak"
$ws.Range("D1").Value = "Quantity"
$ws.Range("E1").Value = "Part Number"
$ws.Range("F1").Value = "No Scan"
$ws.Range("G1").Value = "Tgl Checkin"

# --- Row 2 ---
$ws.Range("A2").Value = "B113"
$ws.Range("B2").Value = "Besar"
$ws.Range("C2").Value = "Penuh"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "4111-03550-C"
$ws.Range("F2").Value = "GQG9DB0W0T4XXP4Q"
$ws.Range("G2").Value = "2023-12-08 15:53:51"

# --- Row 3 ---
$ws.Range("A3").Value = "B113"
$ws.Range("B3").Value = "Besar"
$ws.Range("C3").Value = "Penuh"
$ws.Range("D3").Value = 1000
$ws.Range("E3").Value = "4111-03550-C"
$ws.Range("F3").Value = "GQG9DB0W0T4XXP4Q7"
$ws.Range("G3").Value = "2023-12-11 07:12:19"

# --- Row 4 ---
$ws.Range("A4").Value = "C001"
$ws.Range("B4").Value = "Over Area"
$ws.Range("C4").Value = "Terisi"
$ws.Range("D4").Value = 1000
$ws.Range("E4").Value = "4111-03550-C"
$ws.Range("F4").Value = "GQG9DB0W0T4XXP4Q8"
$ws.Range("G4").Value = "2023-12-11 09:49:37"
